$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 10) ---
$ws.Range("I10").Value = "dias/Costo de compra"
$ws.Range("J10").Value = "costo total de la venta"
$ws.Range("K10").Value = "Dias promedio Invenario"
$ws.Rows.Item(10).RowHeight = 47.25

# --- Column I: dias (simple integers 10..39) ---
for ($r = 11; $r -le 40; $r++) {
    $ws.Cells.Item($r, 9).Value = $r - 1
}

# --- Column J: costo total de la venta = I*F ---
$ws.Range("J11").Formula = "=I11*F11"
$ws.Range("J12:J40").Formula = "=I12*F12"

# --- Column K: dias promedio inventario = H/J ---
$ws.Range("K11").Formula = "=H11/J11"
$ws.Range("K12:K40").Formula = "=H12/J12"

# --- Column L: 1/K, only rows 11-15 ---
$ws.Range("L11").Formula = "=1/K11"
$ws.Range("L12:L15").Formula = "=1/K12"

# M15: average of L11:L15
$ws.Range("M15").Formula = "=SUM(L11:L15)/5"

# Row 16 summary cells (no style)
$ws.Range("L16").Formula = "=SUM(K11:K15)/5"
$ws.Range("M16").Formula = "=5/L16"

# --- Number formats & fill highlighting ---
$ws.Range("K11:K15").NumberFormat = "0.00000"
$ws.Range("K11:K16").Interior.Color = 65535
$ws.Range("K16").NumberFormat = "General"
$ws.Range("K16").Interior.ColorIndex = 0

$ws.Range("L11:L15").NumberFormat = "0.000000"

# --- Totals row 41 ---
$ws.Range("J41").Formula = "=SUM(J11:J40)"
$ws.Range("K41").Formula = "=SUM(K11:K40)"

# --- Row 42 addition ---
$ws.Range("J42").Formula = "=J41/30"

# --- Row 44 (new) ---
$ws.Range("H44").Formula = "=(2000*30)/J42"

# --- Column L width ---
$ws.Columns.Item(12).ColumnWidth = 11.375

# --- Sheet view ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("J42").Select()
